$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2950.5386
$ws.Range("I29").Value = 562.1667
$ws.Range("J29").Value = 4997.7144
$ws.Range("K29").Value = 1686.5001
$ws.Range("L29").Value = 14993.1432
$ws.Range("M29").Value = -1405.5001
$ws.Range("N29").Value = -15555.1432

$ws.Range("H40").Value = 2009697.8
$ws.Range("I40").Value = 2009697.8
$ws.Range("K40").Value = 2009697.8
$ws.Range("M40").Value = -2009522.8

$ws.Range("H51").Value = 4999.5454
$ws.Range("J51").Value = 4999.5454
$ws.Range("L51").Value = 4999.5454
$ws.Range("N51").Value = -5967.5454

$ws.Range("H74").Value = 5304.231
$ws.Range("I74").Value = 2985
$ws.Range("K74").Value = 2985
$ws.Range("M74").Value = -2049

$ws.Range("H77").Value = 5304.231
$ws.Range("I77").Value = 2985
$ws.Range("K77").Value = 14925
$ws.Range("M77").Value = -10245

$ws.Range("H82").Value = 3463.8
$ws.Range("I82").Value = 3463.8
$ws.Range("K82").Value = 10391.4
$ws.Range("M82").Value = -9985.400000000001

$ws.Range("H85").Value = 3463.8
$ws.Range("I85").Value = 3463.8
$ws.Range("K85").Value = 10391.4
$ws.Range("M85").Value = -8987.400000000001

$ws.Range("H98").Value = 1323.5
$ws.Range("I98").Value = 1226.8572
$ws.Range("K98").Value = 1226.8572
$ws.Range("M98").Value = 271.1428000000001

$ws.Range("H122").Value = 1323.5
$ws.Range("I122").Value = 1226.8572
$ws.Range("K122").Value = 3680.5716
$ws.Range("M122").Value = -1230.5716

$ws.Range("H125").Value = 2198.6155
$ws.Range("I125").Value = 1333.3334
$ws.Range("K125").Value = 12000.0006
$ws.Range("M125").Value = -9540.000599999999

$ws.Range("H132").Value = 23811320
$ws.Range("I132").Value = 25001584
$ws.Range("K132").Value = 75004752
$ws.Range("M132").Value = -75002222

$ws.Range("H135").Value = 841.7059
$ws.Range("I135").Value = 748.7692
$ws.Range("K135").Value = 6738.922799999999
$ws.Range("M135").Value = -4203.922799999999

$ws.Range("H137").Value = 2980.1614
$ws.Range("I137").Value = 2735.6086
$ws.Range("K137").Value = 8206.825800000001
$ws.Range("M137").Value = -5656.825800000001

$ws.Range("H138").Value = 2341.8125
$ws.Range("I138").Value = 1060.5
$ws.Range("J138").Value = 3110.6
$ws.Range("K138").Value = 3181.5
$ws.Range("L138").Value = 9331.799999999999
$ws.Range("M138").Value = 1958.5
$ws.Range("N138").Value = -19611.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9069.218999999999
$ws.Range("I32").Value = 8281.861999999999
$ws.Range("K32").Value = 8281.861999999999
$ws.Range("M32").Value = -7994.861999999999

$ws.Range("H61").Value = 20837420
$ws.Range("I61").Value = 25003354
$ws.Range("K61").Value = 25003354
$ws.Range("M61").Value = -25003142

$ws.Range("H132").Value = 23289862
$ws.Range("I132").Value = 1656.6571
$ws.Range("J132").Value = 125175770
$ws.Range("K132").Value = 4969.971299999999
$ws.Range("L132").Value = 375527310
$ws.Range("M132").Value = -2439.971299999999
$ws.Range("N132").Value = -375532370

$ws.Range("H136").Value = 20837420
$ws.Range("I136").Value = 25003354
$ws.Range("K136").Value = 75010062
$ws.Range("M136").Value = -75007512

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2405.5483
$ws.Range("I134").Value = 2185.75
$ws.Range("K134").Value = 6557.25
$ws.Range("M134").Value = -4022.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 196.92592
$ws.Range("I7").Value = 46.25
$ws.Range("J7").Value = 416.0909
$ws.Range("K7").Value = 46.25
$ws.Range("L7").Value = 416.0909
$ws.Range("M7").Value = 66.75
$ws.Range("N7").Value = -642.0908999999999

$ws.Range("H31").Value = 5852.7646
$ws.Range("I31").Value = 3457
$ws.Range("K31").Value = 3457
$ws.Range("M31").Value = -3162

$ws.Range("H34").Value = 5852.7646
$ws.Range("I34").Value = 3457
$ws.Range("K34").Value = 3457
$ws.Range("M34").Value = -3255

$ws.Range("H94").Value = 2170.9167
$ws.Range("J94").Value = 2215.15
$ws.Range("L94").Value = 2215.15
$ws.Range("N94").Value = -3117.15

$ws.Range("H99").Value = 24175.916
$ws.Range("I99").Value = 25708.818
$ws.Range("J99").Value = 7314
$ws.Range("K99").Value = 25708.818
$ws.Range("L99").Value = 7314
$ws.Range("M99").Value = -24210.818
$ws.Range("N99").Value = -10310

$ws.Range("H126").Value = 24175.916
$ws.Range("I126").Value = 25708.818
$ws.Range("J126").Value = 7314
$ws.Range("K126").Value = 77126.454
$ws.Range("L126").Value = 21942
$ws.Range("M126").Value = -74656.454
$ws.Range("N126").Value = -26882

$ws.Range("H134").Value = 1777.5
$ws.Range("I134").Value = 1633.6
$ws.Range("K134").Value = 4900.799999999999
$ws.Range("M134").Value = -2365.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7632
$ws.Range("J3").Value = 7995
$ws.Range("L3").Value = 23985
$ws.Range("N3").Value = -24209

$ws.Range("H55").Value = 6667086.5
$ws.Range("I55").Value = 442.7143
$ws.Range("J55").Value = 22222588
$ws.Range("K55").Value = 1328.1429
$ws.Range("L55").Value = 66667764
$ws.Range("M55").Value = -1151.1429
$ws.Range("N55").Value = -66668118

$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868

$ws.Range("H107").Value = 1225.3889
$ws.Range("I107").Value = 611.1667
$ws.Range("J107").Value = 1532.5
$ws.Range("K107").Value = 1833.5001
$ws.Range("L107").Value = 4597.5
$ws.Range("M107").Value = 86.49990000000003
$ws.Range("N107").Value = -8437.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3515
$ws.Range("I102").Value = 2089.1428
$ws.Range("K102").Value = 2089.1428
$ws.Range("M102").Value = -467.1428000000001

$ws.Range("H122").Value = 4287.5
$ws.Range("I122").Value = 4116.6665
$ws.Range("K122").Value = 12349.9995
$ws.Range("M122").Value = -9899.999500000002

$ws.Range("H126").Value = 4026.5
$ws.Range("I126").Value = 3166
$ws.Range("J126").Value = 4542.8
$ws.Range("K126").Value = 9498
$ws.Range("L126").Value = 13628.4
$ws.Range("N126").Value = -18568.4
$ws.Range("M126").Value = -7028

$ws.Range("H132").Value = 1330.9
$ws.Range("I132").Value = 1214.75
$ws.Range("J132").Value = 1795.5
$ws.Range("K132").Value = 3644.25
$ws.Range("L132").Value = 5386.5
$ws.Range("M132").Value = -1114.25
$ws.Range("N132").Value = -10446.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4024
$ws.Range("I7").Value = 3580.8333
$ws.Range("J7").Value = 4289.9
$ws.Range("K7").Value = 3580.8333
$ws.Range("L7").Value = 4289.9
$ws.Range("M7").Value = -3468.8333
$ws.Range("N7").Value = -4513.9

$ws.Range("H22").Value = 3708.3333
$ws.Range("I22").Value = 1332.6666
$ws.Range("J22").Value = 6084
$ws.Range("K22").Value = 1332.6666
$ws.Range("L22").Value = 6084
$ws.Range("M22").Value = -1037.6666
$ws.Range("N22").Value = -6674

$ws.Range("H27").Value = 3708.3333
$ws.Range("I27").Value = 1332.6666
$ws.Range("J27").Value = 6084
$ws.Range("K27").Value = 1332.6666
$ws.Range("L27").Value = 6084
$ws.Range("M27").Value = -1225.6666
$ws.Range("N27").Value = -6298

$ws.Range("H61").Value = 2308.9644
$ws.Range("I61").Value = 1669.6428
$ws.Range("K61").Value = 1669.6428
$ws.Range("M61").Value = -1467.6428

$ws.Range("H113").Value = 2308.9644
$ws.Range("I113").Value = 1669.6428
$ws.Range("K113").Value = 1669.6428
$ws.Range("M113").Value = 500.3571999999999

$ws.Range("H126").Value = 4024
$ws.Range("I126").Value = 3580.8333
$ws.Range("J126").Value = 4289.9
$ws.Range("K126").Value = 10742.4999
$ws.Range("L126").Value = 12869.7
$ws.Range("M126").Value = -8272.499899999999
$ws.Range("N126").Value = -17809.7

$ws.Range("H132").Value = 19012.678
$ws.Range("I132").Value = 21470.545
$ws.Range("J132").Value = 17422.295
$ws.Range("K132").Value = 64411.63499999999
$ws.Range("L132").Value = 52266.88499999999
$ws.Range("M132").Value = -61881.63499999999
$ws.Range("N132").Value = -57326.88499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 127587.25
$ws.Range("I122").Value = 127587.25
$ws.Range("K122").Value = 382761.75
$ws.Range("M122").Value = -380311.75

$ws.Range("H126").Value = 12159.8
$ws.Range("I126").Value = 12159.8
$ws.Range("K126").Value = 36479.39999999999
$ws.Range("M126").Value = -34009.39999999999

$ws.Range("H132").Value = 4097.8125
$ws.Range("I132").Value = 4304.3335
$ws.Range("K132").Value = 12913.0005
$ws.Range("M132").Value = -10383.0005

$ws.Range("H133").Value = 75750
$ws.Range("J133").Value = 75750
$ws.Range("L133").Value = 75750
$ws.Range("N133").Value = -85870
